# (#91) Added Tests for Basic Search form
$wb = $excel.ActiveWorkbook

# --- cts_pages (sheet 1): move selection from C3 to A3 ---
$ws1 = $wb.Worksheets.Item("cts_pages")
$ws1.Range("A3").Select() | Out-Null

# --- add the new "basicsearch_form" sheet after the last existing sheet ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "basicsearch_form"

# Fill data column by column (matches the shared-string order of the source workbook)
$ws3.Range("A1").Value = "path"
$ws3.Range("A2").Value = "/about-cancer/treatment/clinical-trials/search/"

$ws3.Range("B1").Value = "header"
$ws3.Range("B2").Value = "Find NCI-Supported Clinical Trials"

$ws3.Range("C1").Value = "lbl_cancertype"
$ws3.Range("C2").Value = "Cancer Type/Keyword"

$ws3.Range("D1").Value = "helptext_cancertype"

$ws3.Range("E1").Value = "helpicon_cancertype"

$ws3.Range("F1").Value = "lbl_age"
$ws3.Range("F2").Value = "Cancer Type/Keyword"

$ws3.Range("G1").Value = "helptext_age"

$ws3.Range("H1").Value = "helpicon_age"

# Column widths (best fit on header/value columns)
$ws3.Columns.Item(1).ColumnWidth = 39.1
$ws3.Columns.Item(3).ColumnWidth = 11.75
$ws3.Columns.Item(4).ColumnWidth = 16.75
$ws3.Columns.Item(6).ColumnWidth = 11.75
$ws3.Columns.Item(7).ColumnWidth = 16.75

# Selection on the new sheet
$ws3.Range("G1").Select() | Out-Null
